$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1826.6666
$ws.Range("I40").Value = 1500
$ws.Range("J40").Value = 1850
$ws.Range("K40").Value = 1500
$ws.Range("L40").Value = 1850
$ws.Range("M40").Value = -1325
$ws.Range("N40").Value = -2200
$ws.Range("H64").Value = 2201460.8
$ws.Range("I64").Value = 3348073.2
$ws.Range("J64").Value = 3786.6667
$ws.Range("K64").Value = 3348073.2
$ws.Range("L64").Value = 3786.6667
$ws.Range("M64").Value = -3347825.2
$ws.Range("N64").Value = -4282.6667
$ws.Range("H67").Value = 2201460.8
$ws.Range("I67").Value = 3348073.2
$ws.Range("J67").Value = 3786.6667
$ws.Range("K67").Value = 3348073.2
$ws.Range("L67").Value = 3786.6667
$ws.Range("M67").Value = -3347215.2
$ws.Range("N67").Value = -5502.6667
$ws.Range("H76").Value = 3900.625
$ws.Range("I76").Value = 3882.1875
$ws.Range("J76").Value = 3937.5
$ws.Range("K76").Value = 3882.1875
$ws.Range("L76").Value = 3937.5
$ws.Range("M76").Value = -3567.1875
$ws.Range("N76").Value = -4567.5
$ws.Range("H79").Value = 3900.625
$ws.Range("I79").Value = 3882.1875
$ws.Range("J79").Value = 3937.5
$ws.Range("K79").Value = 3882.1875
$ws.Range("L79").Value = 3937.5
$ws.Range("M79").Value = -2790.1875
$ws.Range("N79").Value = -6121.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 62503236
$ws.Range("I63").Value = 100002580
$ws.Range("J63").Value = 4334.6665
$ws.Range("K63").Value = 100002580
$ws.Range("L63").Value = 4334.6665
$ws.Range("M63").Value = -100001894
$ws.Range("N63").Value = -5706.6665
$ws.Range("H66").Value = 62503236
$ws.Range("I66").Value = 100002580
$ws.Range("J66").Value = 4334.6665
$ws.Range("K66").Value = 500012900
$ws.Range("L66").Value = 21673.3325
$ws.Range("M66").Value = -500009468
$ws.Range("N66").Value = -28537.3325
$ws.Range("H88").Value = 1884.4286
$ws.Range("J88").Value = 1998.2
$ws.Range("L88").Value = 1998.2
$ws.Range("N88").Value = -2810.2
$ws.Range("H91").Value = 1884.4286
$ws.Range("J91").Value = 1998.2
$ws.Range("L91").Value = 1998.2
$ws.Range("N91").Value = -4806.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 23812152
$ws.Range("I86").Value = 2427.0667
$ws.Range("J86").Value = 83336460
$ws.Range("K86").Value = 2427.0667
$ws.Range("L86").Value = 83336460
$ws.Range("M86").Value = -1304.0667
$ws.Range("N86").Value = -83338706
$ws.Range("H89").Value = 23812152
$ws.Range("I89").Value = 2427.0667
$ws.Range("J89").Value = 83336460
$ws.Range("K89").Value = 12135.3335
$ws.Range("L89").Value = 416682300
$ws.Range("M89").Value = -6519.333499999999
$ws.Range("N89").Value = -416693532
$ws.Range("H94").Value = 640.8484999999999
$ws.Range("I94").Value = 467.9565
$ws.Range("J94").Value = 1038.5
$ws.Range("K94").Value = 467.9565
$ws.Range("L94").Value = 1038.5
$ws.Range("M94").Value = -16.95650000000001
$ws.Range("N94").Value = -1940.5
$ws.Range("H99").Value = 874.25
$ws.Range("I99").Value = 855.8889
$ws.Range("J99").Value = 929.3333
$ws.Range("K99").Value = 855.8889
$ws.Range("L99").Value = 929.3333
$ws.Range("M99").Value = 642.1111
$ws.Range("N99").Value = -3925.3333
$ws.Range("H105").Value = 3976.5217
$ws.Range("I105").Value = 3271
$ws.Range("J105").Value = 4519.231
$ws.Range("K105").Value = 3271
$ws.Range("L105").Value = 4519.231
$ws.Range("M105").Value = -1524
$ws.Range("N105").Value = -8013.231

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2374.72
$ws.Range("I62").Value = 2250
$ws.Range("J62").Value = 2509.8333
$ws.Range("K62").Value = 2250
$ws.Range("L62").Value = 2509.8333
$ws.Range("M62").Value = -1626
$ws.Range("N62").Value = -3757.8333
$ws.Range("H65").Value = 2374.72
$ws.Range("I65").Value = 2250
$ws.Range("J65").Value = 2509.8333
$ws.Range("K65").Value = 11250
$ws.Range("L65").Value = 12549.1665
$ws.Range("M65").Value = -8130
$ws.Range("N65").Value = -18789.1665
$ws.Range("H86").Value = 3286.1853
$ws.Range("I86").Value = 3044.2144
$ws.Range("J86").Value = 3546.7693
$ws.Range("K86").Value = 3044.2144
$ws.Range("L86").Value = 3546.7693
$ws.Range("M86").Value = -1921.2144
$ws.Range("N86").Value = -5792.7693
$ws.Range("H89").Value = 3286.1853
$ws.Range("I89").Value = 3044.2144
$ws.Range("J89").Value = 3546.7693
$ws.Range("K89").Value = 15221.072
$ws.Range("L89").Value = 17733.8465
$ws.Range("M89").Value = -9605.072
$ws.Range("N89").Value = -28965.8465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 4777
$ws.Range("I3").Value = 2618.6
$ws.Range("J3").Value = 7475
$ws.Range("K3").Value = 7855.799999999999
$ws.Range("L3").Value = 22425
$ws.Range("M3").Value = -7743.799999999999
$ws.Range("N3").Value = -22649

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 16180.9375
$ws.Range("I70").Value = 200000
$ws.Range("J70").Value = 3926.3333
$ws.Range("K70").Value = 200000
$ws.Range("L70").Value = 3926.3333
$ws.Range("M70").Value = -199730
$ws.Range("N70").Value = -4466.3333
$ws.Range("H73").Value = 16180.9375
$ws.Range("I73").Value = 200000
$ws.Range("J73").Value = 3926.3333
$ws.Range("K73").Value = 200000
$ws.Range("L73").Value = 3926.3333
$ws.Range("M73").Value = -199064
$ws.Range("N73").Value = -5798.3333
$ws.Range("H80").Value = 13891643
$ws.Range("I80").Value = 25643388
$ws.Range("J80").Value = 3218
$ws.Range("K80").Value = 25643388
$ws.Range("L80").Value = 3218
$ws.Range("M80").Value = -25642390
$ws.Range("N80").Value = -5214
$ws.Range("H83").Value = 13891643
$ws.Range("I83").Value = 25643388
$ws.Range("J83").Value = 3218
$ws.Range("K83").Value = 128216940
$ws.Range("L83").Value = 16090
$ws.Range("M83").Value = -128211948
$ws.Range("N83").Value = -26074

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1747.1666
$ws.Range("I82").Value = 1596.6
$ws.Range("J82").Value = 2500
$ws.Range("K82").Value = 1596.6
$ws.Range("L82").Value = 2500
$ws.Range("M82").Value = -1235.6
$ws.Range("N82").Value = -3222
$ws.Range("H85").Value = 1747.1666
$ws.Range("I85").Value = 1596.6
$ws.Range("J85").Value = 2500
$ws.Range("K85").Value = 1596.6
$ws.Range("L85").Value = 2500
$ws.Range("M85").Value = -348.5999999999999
$ws.Range("N85").Value = -4996
